$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.565.06"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "3.177.07"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'599.70"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").Value = "'152.07"
$ws.Range("E6").Value = "  -3.82%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.178.45"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -3.52%  "
$ws.Range("E10").Value = "  -5.34%  "
$ws.Range("D11").Value = "'5.52"
$ws.Range("E11").Value = "  -3.21%  "
$ws.Range("D12").Value = "'0.475"
$ws.Range("E12").Value = "  -6.15%  "
$ws.Range("E13").Value = "  -5.74%  "
$ws.Range("D14").Value = "'36.93"
$ws.Range("E14").Value = "  -5.38%  "
$ws.Range("D15").Value = "3.705.52"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").Value = "64.643.21"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "3.190.99"
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "'7.03"
$ws.Range("E19").Value = "  -5.02%  "
$ws.Range("D20").Value = "'480.32"
$ws.Range("E20").Value = "  -5.91%  "
$ws.Range("D21").Value = "'14.81"
$ws.Range("E21").Value = "  -3.24%  "
$ws.Range("D22").Value = "'0.715"
$ws.Range("E22").Value = "  -2.83%  "
$ws.Range("D23").Value = "'7.75"
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("D24").Value = "'13.88"
$ws.Range("E24").Value = "  -5.90%  "
$ws.Range("D25").Value = "'84.76"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  -0.26%  "
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").Value = "'8.62"
$ws.Range("E29").Value = "  -5.73%  "
$ws.Range("D30").Value = "'0.122"
$ws.Range("E30").Value = "  +15.71%  "
$ws.Range("D31").Value = "'6.96"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("E32").Value = "  -8.59%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'26.86"
$ws.Range("E34").Value = "  -4.87%  "
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("D36").Value = "'6.12"
$ws.Range("E36").Value = "  -5.99%  "
$ws.Range("D37").Value = "'54.65"
$ws.Range("E37").Value = "  -2.51%  "
$ws.Range("D38").Value = "'3.25"
$ws.Range("E38").Value = "  +6.43%  "
$ws.Range("E39").Value = "  -4.02%  "
$ws.Range("D40").Value = "'457.31"
$ws.Range("E40").Value = "  -10.75%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.125"
$ws.Range("E41").Value = "  -3.93%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0402"
$ws.Range("E42").Value = "  -4.76%  "
$ws.Range("D43").Value = "'8.51"
$ws.Range("E43").Value = "  -2.99%  "
$ws.Range("D44").Value = "'2.42"
$ws.Range("E44").Value = "  -2.35%  "
$ws.Range("D45").Value = "2.876.11"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").Value = "'0.275"
$ws.Range("E46").Value = "  -8.11%  "
$ws.Range("D47").Value = "'27.34"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("D49").Value = "'2.34"
$ws.Range("E49").Value = "  -3.04%  "
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "'120.19"
$ws.Range("E51").Value = "  -1.84%  "
